$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("B2").Value = 73.7001285719821
$ws.Range("C2").Value = 72.3336424210523
$ws.Range("D2").Value = 75.0666147229119
$ws.Range("B3").Value = 76.239726224251
$ws.Range("C3").Value = 74.2016212947435
$ws.Range("D3").Value = 78.2778311537585
$ws.Range("B4").Value = 79.2192155789382
$ws.Range("C4").Value = 74.3248868209307
$ws.Range("D4").Value = 84.1135443369457
$ws.Range("B5").Value = 76.5089904149901
$ws.Range("C5").Value = 72.0562147796462
$ws.Range("D5").Value = 80.9617660503339
$ws.Range("B10").Value = 70.6179655561316
$ws.Range("C10").Value = 63.4813668790681
$ws.Range("D10").Value = 77.754564233195
$ws.Range("B11").Value = 71.9315392433974
$ws.Range("C11").Value = 68.5573648738596
$ws.Range("D11").Value = 75.3057136129352
$ws.Range("B12").Value = 77.8385029514733
$ws.Range("C12").Value = 73.4720101468262
$ws.Range("D12").Value = 82.2049957561204
$ws.Range("B14").Value = 65.6756128159912
$ws.Range("C14").Value = 62.6818047274839
$ws.Range("D14").Value = 68.6694209044986
$ws.Range("B15").Value = 68.7293820559067
$ws.Range("C15").Value = 67.2928692548832
$ws.Range("D15").Value = 70.1658948569303
$ws.Range("B16").Value = 70.0349515464357
$ws.Range("C16").Value = 67.8340958209815
$ws.Range("D16").Value = 72.23580727189
$ws.Range("B17").Value = 73.8011174594306
$ws.Range("C17").Value = 68.7699381998794
$ws.Range("D17").Value = 78.8322967189819
$ws.Range("B18").Value = 71.9545300671651
$ws.Range("C18").Value = 67.3409803626745
$ws.Range("D18").Value = 76.5680797716556
$ws.Range("B23").Value = 57.093643885943
$ws.Range("C23").Value = 49.1509715456626
$ws.Range("D23").Value = 65.0363162262233
$ws.Range("B24").Value = 63.7863454580277
$ws.Range("C24").Value = 60.1526597008331
$ws.Range("D24").Value = 67.4200312152222
$ws.Range("B25").Value = 69.8931601175901
$ws.Range("C25").Value = 65.0588305973776
$ws.Range("D25").Value = 74.7274896378027
$ws.Range("B27").Value = 64.7875368806995
$ws.Range("C27").Value = 61.8104609935632
$ws.Range("D27").Value = 67.7646127678359
$ws.Range("B28").Value = 67.4426950982622
$ws.Range("C28").Value = 65.9938812820719
$ws.Range("D28").Value = 68.8915089144525
$ws.Range("B29").Value = 69.5524559847153
$ws.Range("C29").Value = 67.3382475538882
$ws.Range("D29").Value = 71.7666644155423
$ws.Range("B30").Value = 69.2641485591769
$ws.Range("C30").Value = 63.4068241328164
$ws.Range("D30").Value = 75.1214729855374
$ws.Range("B31").Value = 68.5089970578903
$ws.Range("C31").Value = 63.5248877895078
$ws.Range("D31").Value = 73.4931063262728
$ws.Range("B36").Value = 63.9965718289014
$ws.Range("C36").Value = 56.5557068223802
$ws.Range("D36").Value = 71.4374368354227
$ws.Range("B37").Value = 58.6717978225169
$ws.Range("C37").Value = 54.8678003845487
$ws.Range("D37").Value = 62.4757952604851
$ws.Range("B38").Value = 74.0232378134785
$ws.Range("C38").Value = 69.1540095101055
$ws.Range("D38").Value = 78.8924661168515
$ws.Range("B40").Value = 62.4616898045869
$ws.Range("C40").Value = 59.5216656437741
$ws.Range("D40").Value = 65.4017139653997
